# Auto update: 2025-12-03 03:05:24
# Updates the DECISION/국장_방산_분석 sheet with refreshed market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a plain-text value into a cell without letting Excel's
# input parser reinterpret a date-shaped string (e.g. "2025-12-03") as
# a real date serial. We build the text via a formula (so the result is
# a literal text value with default/General formatting), copy it, and
# paste-special just the values into the destination - this keeps the
# destination cell's number format/style completely untouched.
# ---------------------------------------------------------------------
$ws.Range("Z1").Formula = "=""2025-12-03"""
$ws.Range("Z1").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("A3").PasteSpecial(-4163)
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Row 2 - KOREA AEROSPACE (047810.KS)
$ws.Range("D2").Value = 105400
$ws.Range("E2").Value = 52.6
$ws.Range("F2").Value = -2.95
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 53
$ws.Range("K2").Value = 55.6
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 65.32892478746797
$ws.Range("O2").Value = "🟢 상승 우위 (다소 완화)"

# Row 3 - HYUNDAI ROTEM (064350.KS)
$ws.Range("D3").Value = 169500
$ws.Range("E3").Value = 20.1
$ws.Range("F3").Value = -4.35
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 50.6
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 65.32892478746797
$ws.Range("O3").Value = "🟢 상승 우위 (다소 완화)"

# Row 4 - HANWHA SYSTEMS (272210.KS)
$ws.Range("D4").Value = 45550
$ws.Range("E4").Value = 14.9
$ws.Range("F4").Value = -1.41
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 70
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 50.6
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 65.32892478746797
$ws.Range("O4").Value = "🟢 상승 우위 (다소 완화)"

# Row 5 - LIG Nex1 (079550.KS)
$ws.Range("D5").Value = 364500
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = -5.2
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 66
$ws.Range("I5").Value = 66
$ws.Range("J5").Value = 63
$ws.Range("K5").Value = 49
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 65.32892478746797
$ws.Range("O5").Value = "🟢 상승 우위 (다소 완화)"

# Row 6 - HANWHA AEROSPACE (012450.KS)
$ws.Range("D6").Value = 824000
$ws.Range("E6").Value = 23.8
$ws.Range("F6").Value = -4.07
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 53
$ws.Range("I6").Value = 63
$ws.Range("J6").Value = 66
$ws.Range("K6").Value = 47.8
$ws.Range("M6").Value = "⛔ 관망하십시오."
$ws.Range("N6").Value = 65.32892478746797
$ws.Range("O6").Value = "🟢 상승 우위 (다소 완화)"
